$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Estado de Cuenta" worksheet lists one row per overdue period
# (2505/2506/2507 before the edit). This update adds a new period
# (2508) to the bottom of the list and shifts the per-period values
# down by one period, so:
#   old row16 (period 2507) data now belongs to row18 (period 2507 stays)
#   old row18 (period 2505, valor mora 56506) data now belongs to row16
#   a brand-new row19 is added for period 2508 (valor mora 65200)
# The signature block below (rows 23/24) shifts down by one row as a
# consequence (now rows 24/25).
# ------------------------------------------------------------------

# Insert a new row before row 19; this pushes the old rows 19-24 (the
# blank rows plus the signature block) down to 20-25.
$ws.Rows("19").Insert()

# The freshly inserted row 19 is blank. Row 18 still carries the
# "last data row" look (bottom border on every cell), so copy that
# formatting down onto the new row 19 first.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 18 is no longer the last data row, so it should look like the
# other "middle" rows (16/17) instead - copy row 16's formatting onto it.
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the period / valor-mora values for each data row ---

# Row 16: period 2505, valor mora 56506 (salario basico stays 1630000)
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 56506

# Row 17: period 2506, valor mora 65200 -> unchanged

# Row 18: period 2507, valor mora 65200
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 65200

# Row 19 (new): period 2508, valor mora 65200
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1052969719"
$ws.Range("D19").Value = "RAUL FERNANDO SIMANCA MARTINEZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 65200
$ws.Range("G19").Value = 1630000

# --- Update the summary fields at the top of the sheet ---

# VALOR MORA total = sum of the four periods' valor mora (56506+65200*3)
$ws.Range("E11").Value = 252106

# Cant. Periodos now counts 4 overdue periods instead of 3
$ws.Range("F13").Value = 4
